# Append a new paragraph after the last existing paragraph ("My name is
# karan.") containing the line "This file is been added to .gitignore.",
# matching the formatting (en-US language run properties) of the
# surrounding text.
$d = $word.ActiveDocument

$end = $d.Content
$end.Collapse(0)            # wdCollapseEnd - collapse to the very end of the doc
$end.InsertParagraphAfter() # insert a new empty paragraph after the current end
$end.Collapse(0)            # move into the newly-created paragraph
$end.MoveEnd(1, 1)          # wdCharacter - extend range onto the new paragraph mark
$end.Collapse(0)            # collapse again so we're positioned inside the new paragraph
$end.Text = "This file is been added to .gitignore."
